# Apply updated NATMI TPM recomputation to the LR-pairs sheet.
# Adds self-cluster (sending == target) pairs and refreshes all metrics
# for the Bsg -> Slc16a7 ligand-receptor table (rows 2-17, all 20 columns).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs", "Bsg", "Slc16a7", "ECs", 3, 1, 49.07461933333334, 147.223858, 0.2082521693470354, 0.2082521693470354, 2, 0.6666666666666666, 0.1713566666666667, 0.51407, 0.01023284555714517, 0.01023284555714517, 8.409263186895558, 75.68336868206, 0.002131012285868655, 0.002131012285868654),
    @("ECs", "Bsg", "Slc16a7", "FAPs", 3, 1, 49.07461933333334, 147.223858, 0.2082521693470354, 0.2082521693470354, 3, 1, 0.6847623333333334, 2.054287, 0.04089171046949069, 0.04089171046949069, 33.60445084213845, 302.440057579246, 0.008515787413582316, 0.008515787413582313),
    @("ECs", "Bsg", "Slc16a7", "MuSCs", 3, 1, 49.07461933333334, 147.223858, 0.2082521693470354, 0.2082521693470354, 3, 1, 0.8339630000000001, 2.501889, 0.04980147399793875, 0.04980147399793875, 40.92641676308467, 368.3377508677621, 0.01037126499675072, 0.01037126499675072),
    @("ECs", "Bsg", "Slc16a7", "Resolving-Mac", 3, 1, 49.07461933333334, 147.223858, 0.2082521693470354, 0.2082521693470354, 3, 1, 15.05566733333333, 45.167002, 0.8990739699754254, 0.8990739699754254, 738.8511431926352, 6649.660288733716, 0.1872341046508337, 0.1872341046508337),
    @("FAPs", "Bsg", "Slc16a7", "ECs", 3, 1, 111.3149186666667, 333.944756, 0.4723739808466805, 0.4723739808466804, 2, 0.6666666666666666, 0.1713566666666667, 0.51407, 0.01023284555714517, 0.01023284555714517, 19.07455341299111, 171.67098071692, 0.004833729991217932, 0.004833729991217931),
    @("FAPs", "Bsg", "Slc16a7", "FAPs", 3, 1, 111.3149186666667, 333.944756, 0.4723739808466805, 0.4723739808466804, 3, 1, 0.6847623333333334, 2.054287, 0.04089171046949069, 0.04089171046949069, 76.22426344099688, 686.0183709689719, 0.0193161800581032, 0.01931618005810319),
    @("FAPs", "Bsg", "Slc16a7", "MuSCs", 3, 1, 111.3149186666667, 333.944756, 0.4723739808466805, 0.4723739808466804, 3, 1, 0.8339630000000001, 2.501889, 0.04980147399793875, 0.04980147399793875, 92.83252351600933, 835.492711644084, 0.02352492052443877, 0.02352492052443877),
    @("FAPs", "Bsg", "Slc16a7", "Resolving-Mac", 3, 1, 111.3149186666667, 333.944756, 0.4723739808466805, 0.4723739808466804, 3, 1, 15.05566733333333, 45.167002, 0.8990739699754254, 0.8990739699754254, 1675.92038468239, 15083.28346214151, 0.4246991502729206, 0.4246991502729205),
    @("MuSCs", "Bsg", "Slc16a7", "ECs", 3, 1, 40.34671033333333, 121.040131, 0.1712145721571795, 0.1712145721571795, 2, 0.6666666666666666, 0.1713566666666667, 0.51407, 0.01023284555714517, 0.01023284555714517, 6.913677793685556, 62.22310014317001, 0.001752012274017105, 0.001752012274017105),
    @("MuSCs", "Bsg", "Slc16a7", "FAPs", 3, 1, 40.34671033333333, 121.040131, 0.1712145721571795, 0.1712145721571795, 3, 1, 0.6847623333333334, 2.054287, 0.04089171046949069, 0.04089171046949069, 27.62790751017745, 248.651167591597, 0.007001256712809106, 0.007001256712809104),
    @("MuSCs", "Bsg", "Slc16a7", "MuSCs", 3, 1, 40.34671033333333, 121.040131, 0.1712145721571795, 0.1712145721571795, 3, 1, 0.8339630000000001, 2.501889, 0.04980147399793875, 0.04980147399793875, 33.64766358971767, 302.828972307459, 0.008526738063353983, 0.008526738063353983),
    @("MuSCs", "Bsg", "Slc16a7", "Resolving-Mac", 3, 1, 40.34671033333333, 121.040131, 0.1712145721571795, 0.1712145721571795, 3, 1, 15.05566733333333, 45.167002, 0.8990739699754254, 0.8990739699754254, 607.4466487730291, 5467.019838957262, 0.1539345651069993, 0.1539345651069993),
    @("Resolving-Mac", "Bsg", "Slc16a7", "ECs", 3, 1, 34.91373066666667, 104.741192, 0.1481592776491046, 0.1481592776491046, 2, 0.6666666666666666, 0.1713566666666667, 0.51407, 0.01023284555714517, 0.01023284555714517, 5.982700507937779, 53.84430457144001, 0.001516091006041478, 0.001516091006041477),
    @("Resolving-Mac", "Bsg", "Slc16a7", "FAPs", 3, 1, 34.91373066666667, 104.741192, 0.1481592776491046, 0.1481592776491046, 3, 1, 0.6847623333333334, 2.054287, 0.04089171046949069, 0.04089171046949069, 23.90760767667823, 215.168469090104, 0.006058486284996068, 0.006058486284996068),
    @("Resolving-Mac", "Bsg", "Slc16a7", "MuSCs", 3, 1, 34.91373066666667, 104.741192, 0.1481592776491046, 0.1481592776491046, 3, 1, 0.8339630000000001, 2.501889, 0.04980147399793875, 0.04980147399793875, 29.11675956796534, 262.0508361116881, 0.007378550413395271, 0.00737855041339527),
    @("Resolving-Mac", "Bsg", "Slc16a7", "Resolving-Mac", 3, 1, 34.91373066666667, 104.741192, 0.1481592776491046, 0.1481592776491046, 3, 1, 15.05566733333333, 45.167002, 0.8990739699754254, 0.8990739699754254, 525.6495142829316, 4730.845628546384, 0.1332061499446718, 0.1332061499446718)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $rowData = $rows[$i]
    $excelRow = $i + 2
    for ($col = 1; $col -le $rowData.Length; $col++) {
        $ws.Cells.Item($excelRow, $col).Value = $rowData[$col - 1]
    }
}

"Updated rows 2-17 with recomputed TPM values"
